# Apply the update described by the diff:
# - Insert a new row at position 46 (pushing existing rows 46-53 down to 47-54)
# - Populate the new row 46 with the new record data
# - The dimension will automatically grow from A1:R53 to A1:R54

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 46, shifting rows 46:53 down to 47:54
$ws.Rows("46:46").Insert()

# Fill in the new row 46 with the new record
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 44504
$ws.Cells.Item(46, 4).NumberFormat = $ws.Cells.Item(47, 4).NumberFormat
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = 100112021
$ws.Cells.Item(46, 7).Value = "Ají"
$ws.Cells.Item(46, 8).Value = "Inferno"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 40
$ws.Cells.Item(46, 11).Value = 14000
$ws.Cells.Item(46, 12).Value = 15000
$ws.Cells.Item(46, 13).Value = 14500
$ws.Cells.Item(46, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 1208
$ws.Cells.Item(46, 17).Value = 12
$ws.Cells.Item(46, 18).Value = "Hortaliza"
